# Updated cryptos list on Tue Aug  6 09:57:38 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to keep the exact text (avoids Excel auto-converting
    # numeric-looking strings like "480.82" or "1.00" into floating point
    # numbers and losing formatting / precision), then restore the
    # cell's original (unstyled) look so no stray formatting is left behind.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "54.865.35"
$ws.Range("E2").Value = "  +5.41%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "2.434.08"
$ws.Range("E3").Value = "  +6.62%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.15%  "

# Row 5 (BNB)
Set-TextValue "D5" "480.82"
$ws.Range("E5").Value = "  +8.92%  "

# Row 6 (Solana)
Set-TextValue "D6" "138.44"
$ws.Range("E6").Value = "  +16.59%  "

# Row 7 (USDC)
Set-TextValue "D7" "0.996"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 (XRP)
$ws.Range("E8").Value = "  +8.60%  "

# Row 9 (LidoStakedEther)
Set-TextValue "D9" "2.449.99"
$ws.Range("E9").Value = "  +7.28%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  +12.52%  "

# Row 11 (Toncoin)
$ws.Range("E11").Value = "  +4.23%  "

# Row 12 (Cardano)
$ws.Range("E12").Value = "  +8.20%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  +2.23%  "

# Row 14 (WrappedliquidstakedEther2.0)
Set-TextValue "D14" "2.849.27"
$ws.Range("E14").Value = "  +7.26%  "

# Row 15 (WrappedBTC)
Set-TextValue "D15" "54.971.60"
$ws.Range("E15").Value = "  +5.87%  "

# Row 16 (Avalanche)
$ws.Range("E16").Value = "  +10.30%  "

# Row 18 (WrappedEther)
Set-TextValue "D18" "2.448.45"
$ws.Range("E18").Value = "  +6.62%  "

# Row 19 (Polkadot)
$ws.Range("E19").Value = "  +10.28%  "

# Row 20 (BitcoinCash)
Set-TextValue "D20" "314.41"
$ws.Range("E20").Value = "  +6.53%  "

# Row 21 (Chainlink)
Set-TextValue "D21" "9.79"
$ws.Range("E21").Value = "  +12.40%  "

# Row 22 (Dai)
Set-TextValue "D22" "0.995"
$ws.Range("E22").Value = "  -0.40%  "

# Row 23 (Uniswap)
$ws.Range("E23").Value = "  +10.36%  "

# Row 24 (Litecoin)
Set-TextValue "D24" "57.11"
$ws.Range("E24").Value = "  +7.83%  "

# Rows 25-27 reordered: Kaspa, Polygon, Binance-PegBSC-USD -> Polygon, Binance-PegBSC-USD, Kaspa
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D25" "0.404"
$ws.Range("E25").Value = "  +11.31%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D27" "0.163"
$ws.Range("E27").Value = "  +15.07%  "

# Row 28 (WrappedeETH)
Set-TextValue "D28" "2.543.98"
$ws.Range("E28").Value = "  +6.38%  "

# Row 29 (InternetComputer(DFINITY))
Set-TextValue "D29" "7.33"
$ws.Range("E29").Value = "  +7.20%  "

# Row 30 (PEPE)
Set-TextValue "D30" "0.0₃0775"
$ws.Range("E30").Value = "  +19.45%  "

# Row 31 (USDe)
Set-TextValue "D31" "0.997"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32 (Monero)
Set-TextValue "D32" "148.69"
$ws.Range("E32").Value = "  +4.39%  "

# Row 33 (EthereumClassic)
Set-TextValue "D33" "17.91"
$ws.Range("E33").Value = "  +7.75%  "

# Row 34 (PancakeSwap)
$ws.Range("E34").Value = "  +12.19%  "

# Row 35 (Aptos)
Set-TextValue "D35" "5.13"
$ws.Range("E35").Value = "  +10.57%  "

# Row 36 (ImmutableX)
$ws.Range("E36").Value = "  +13.75%  "

# Row 37 (NEARProtocol)
$ws.Range("E37").Value = "  +7.73%  "

# Row 38 (Fetch.AI)
Set-TextValue "D38" "0.848"
$ws.Range("E38").Value = "  +5.13%  "

# Rows 39-40 reordered: FirstDigitalUSD, OKB -> OKB, FirstDigitalUSD
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D39" "33.16"
$ws.Range("E39").Value = "  +4.56%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D40" "0.991"
$ws.Range("E40").Value = "  -0.22%  "

# Row 41 (Filecoin)
$ws.Range("E41").Value = "  +10.49%  "

# Rows 42-43 reordered: Hedera, Mantle -> Mantle, Hedera
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D42" "0.597"
$ws.Range("E42").Value = "  +8.58%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D43" "0.0542"
$ws.Range("E43").Value = "  +8.03%  "

# Row 44 (Stacks)
$ws.Range("E44").Value = "  +12.48%  "

# Rows 45-46 reordered: WhiteBITCoin, RenderToken -> RenderToken, WhiteBITCoin
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D45" "4.68"
$ws.Range("E45").Value = "  +14.40%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D46" "10.09"
$ws.Range("E46").Value = "  -0.33%  "

# Row 47 (Bittensor)
Set-TextValue "D47" "253.97"
$ws.Range("E47").Value = "  +29.65%  "

# Row 48 (Stellar)
Set-TextValue "D48" "0.0900"
$ws.Range("E48").Value = "  +10.85%  "

# Row 49 (Maker)
Set-TextValue "D49" "1.935.27"
$ws.Range("E49").Value = "  +2.12%  "

# Row 50 (VeChain)
$ws.Range("E50").Value = "  +9.77%  "

# Row 51 (EnergySwap)
Set-TextValue "D51" "16.97"
$ws.Range("E51").Value = "  +10.48%  "
